# This script re-shuffles the data rows (2-33) of the active worksheet so that
# each target row ends up containing the full set of field values that used to
# live in a different ("source") row, according to the mapping below.
# (targetRow = sourceRow), i.e. row 2 gets the old content of row 6, etc.
# This reproduces a "records got renumbered/reordered" style edit where every
# field of a row moves together (Id, coordinates, species, comments, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (1-based worksheet row numbers)
$mapping = @{
    2  = 6
    3  = 24
    4  = 7
    5  = 13
    6  = 3
    7  = 28
    8  = 4
    9  = 18
    10 = 12
    11 = 8
    12 = 16
    13 = 32
    14 = 11
    15 = 20
    16 = 21
    17 = 10
    18 = 14
    19 = 22
    20 = 31
    21 = 9
    22 = 17
    23 = 2
    24 = 26
    25 = 30
    26 = 27
    27 = 23
    28 = 19
    29 = 25
    30 = 5
    31 = 15
    32 = 33
    33 = 29
}

$firstRow = 2
$lastRow = 33
$lastCol = 51   # column AY

# A few text (inlineStr) columns contain values that look like numbers or dates
# ("1" in column I, "2023-08-28" / "00:00" in columns Y, Z, AA, AB). Force just
# those columns to a text format *before* writing values back, so Excel's
# value-autodetection does not silently convert them to numbers / dates.
$textLookingCols = @(9, 25, 26, 27, 28)   # I, Y, Z, AA, AB

foreach ($c in $textLookingCols) {
    $colRange = $ws.Range($ws.Cells.Item($firstRow, $c), $ws.Cells.Item($lastRow, $c))
    $colRange.NumberFormat = "@"
}

# Read the full block of data (rows 2..33, columns A..AY) into memory once.
# COM returns this as a 2-D array that is 1-based: [1..rowCount, 1..colCount].
$srcRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$original = $srcRange.Value()

# Build the new array applying the row permutation.
# A freshly-created .NET 2-D array is 0-based: [0..rowCount-1, 0..colCount-1].
$rowCount = $lastRow - $firstRow + 1
$new = New-Object 'object[,]' $rowCount, $lastCol

for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $mapping[$targetRow]
    $ti = $targetRow - $firstRow        # 0-based row index into $new
    $si = $sourceRow - $firstRow + 1    # 1-based row index into $original
    for ($c = 1; $c -le $lastCol; $c++) {
        $new[$ti, $c - 1] = $original[$si, $c]
    }
}

# Write the permuted block back in one shot.
$srcRange.Value = $new
